$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- File Handling section: insert two new program rows -------------------
# Current layout (rows 33-39, all under the "File Handling" merged topic):
#   33 File Handling | file_read.py
#   34              | file_write.py
#   35              | file_append.py
#   36              | read_write_binary.py
#   37              | file_unicode.py
#   38              | input.py
#   39              | file_system_cmds.py
#
# Target layout (rows 33-41):
#   33 File Handling | file_read.py
#   34              | file_write.py
#   35              | file_append.py
#   36              | file_overwrite.py   <-- NEW
#   37              | read_write_binary.py
#   38              | file_methods.py     <-- NEW
#   39              | file_unicode.py
#   40              | input.py
#   41              | file_system_cmds.py

# 1) Insert a new row at 36 (everything from 36 down shifts to 37+) and give
#    it the same look (borders/fill/font) as the rest of the File Handling
#    rows by copying the formatting from the row right above it.
$ws.Rows.Item(36).Insert()
$ws.Range("A35:B35").Copy()
$ws.Range("A36:B36").PasteSpecial(-4122)
$ws.Range("A36").Value = ""
$ws.Range("B36").Value = "file_overwrite.py"

# 2) read_write_binary.py is now on row 37. Insert another new row at 38
#    (below it) for file_methods.py, again cloning formatting from a
#    neighboring row in the same block.
$ws.Rows.Item(38).Insert()
$ws.Range("A37:B37").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)
$ws.Range("A38").Value = ""
$ws.Range("B38").Value = "file_methods.py"

$excel.CutCopyMode = $false
